$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A for rows 2-19 from text placeholders to sequential numbers
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# New rows 20-33 data (Column A is numeric id, columns B-H are text)
$newRows = @(
    @(19, "default", "default", "default", "default", "default", "default", "default"),
    @(20, "default", "default", "default", "default", "default", "default", "default"),
    @(21, "sdf", "sdf", "hsdfew", "sdf", "sdf", "sdf", "sdf"),
    @(22, "sdf", "sdf", "hsdfew", "sdf", "sdf", "sdf", "sdf"),
    @(23, "fdh3e", "fdh3e", "dfh24", "fdh3e", "fdh3e", "fdh3e", "fdh3e"),
    @(24, "bfd", "bfd", "gfncvnfd", "bfd", "bfd", "bfd", "bfd"),
    @(25, "hdfw", "hdfw", "rher", "hdfw", "hdfw", "hdfw", "hdfw"),
    @(26, "hdfw", "hdfw", "rher", "hdfw", "hdfw", "hdfw", "hdfw"),
    @(27, "fsdg", "fsdg", "sgsg", "fsdg", "fsdg", "fsdg", "fsdg"),
    @(28, "wert", "wert", "hdf", "wert", "wert", "wert", "wert"),
    @(29, "gsdf", "gfdhjy", "sdfcvsdg", "fdgsdg", "bcvbsdf", "asfvcx", "sfcxv"),
    @(30, "sdfhtkuy", "nvbndfghvbn", "fgdcv", "fgdfbcvb", "fdgfdvc", "bfgdfg", "bdsdgg"),
    @(31, "sdfhtkuy", "nvbndfghvbn", "fgdcv", "fgdfbcvb", "fdgfdvc", "bfgdfg", "bdsdgg"),
    @(32, "sdgcxv", "sdfxc", "vasfsd", "xcvsadf", "vxcasf", "bvfgdfj", "cvsdfs")
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $rowIndex = $startRow + $i
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value = $rowData[$c - 1]
    }
}
